# Applies the WTREGEN.xlsx update:
#  - appends 6 new weekly observation rows (104-109) to the "Data" sheet
#  - refreshes the FRED series metadata dates/popularity on "SeriesInfo"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Data" sheet - append new weekly rows after the existing last row
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$newRows = @(
    @{ Row = 104; Date = 45189; Value = 640.296 },
    @{ Row = 105; Date = 45196; Value = 681.143 },
    @{ Row = 106; Date = 45203; Value = 668.009 },
    @{ Row = 107; Date = 45210; Value = 709.16 },
    @{ Row = 108; Date = 45217; Value = 759.878 },
    @{ Row = 109; Date = 45224; Value = 834.418 }
)

foreach ($r in $newRows) {
    $data.Range("A$($r.Row)").Value = $r.Date
    $data.Range("B$($r.Row)").Value = $r.Value
}

# Copy the date-column formatting (style) from the preceding row (A103)
# onto the newly-added date cells so they match the rest of the column.
$data.Range("A103").Copy()
$data.Range("A104:A109").PasteSpecial(-4122)
$data.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) "SeriesInfo" sheet - update the metadata values
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end: plain date-looking strings. Prefix with
# an apostrophe so Excel stores them as text instead of auto-converting
# to a date serial number, then reset the style back to Normal so no
# stray number-format is left attached to the cell.
$info.Range("B3").Value = "'2023-10-27"
$info.Range("B3").Style = "Normal"

$info.Range("B4").Value = "'2023-10-27"
$info.Range("B4").Style = "Normal"

# observation_end
$info.Range("B7").Value = "'2023-10-25"
$info.Range("B7").Style = "Normal"

# last_updated (includes a time + UTC offset, Excel won't auto-parse it,
# but keep it consistent with the rest regardless).
$info.Range("B14").Value = "2023-10-26 15:34:02-05"

# popularity (numeric)
$info.Range("B15").Value = 78

Write-Output "WTREGEN.xlsx update applied"
